# "Ações com Tabela adicionadas"
# Adds a new worksheet "TC_003_Teste_Tabela_Campo_Teste" (a copy of the
# "TC_002_Teste_Alert_Campo_Teste" sheet, same header/data rows + styles),
# placed right after TC_002, makes it the active/selected sheet, and
# resets the selection on TC_002 back to the full data range.

$wb = $excel.ActiveWorkbook

# Source sheet to duplicate.
$ws2 = $wb.Worksheets.Item("TC_002_Teste_Alert_Campo_Teste")

# Copy it, placing the new sheet immediately after itself.
$ws2.Copy($null, $ws2)

# The newly created copy becomes the sheet right after $ws2.
$ws3 = $wb.Worksheets.Item($ws2.Index + 1)
$ws3.Name = "TC_003_Teste_Tabela_Campo_Teste"

# Reset TC_002's own selection to the whole used range (no longer the
# active/tab-selected sheet).
$ws2.Activate() | Out-Null
$ws2.Range("A1:D2").Select() | Out-Null

# The new sheet becomes the active tab with cell D8 selected.
$ws3.Activate() | Out-Null
$ws3.Range("D8").Select() | Out-Null
